# Lab04: update stack-trace spreadsheet to match assembly code.
# The contents of B78:E100 (minus A, which holds the address ladder
# formulas) shift down by three rows; rows 78-80 end up empty and three
# new rows of data appear at the bottom (98-99 get B values, 100 gets a
# refreshed D/E pair).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear cells that have no content after the edit ---
$clearRefs = @(
    "C78", "C79", "C80",
    "D81", "E81",
    "D82", "E82",
    "B84",
    "D87", "E87",
    "D88", "E88",
    "B89",
    "B91",
    "D99", "E99"
)
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}

# --- set the new cell values (shifted content + new tail rows) ---
$ws.Range("C81").Value = "int product"

$ws.Range("C82").Value = "Save R2"

$ws.Range("C83").Value = "Save R1"

$ws.Range("C84").Value = "previous frame pointer"
$ws.Range("D84").Value = "R5"
$ws.Range("E84").Value = "current frame pointer"

$ws.Range("C85").Value = "square() return address"
$ws.Range("D85").Value = "R7"
$ws.Range("E85").Value = "current return address"

$ws.Range("C86").Value = "int x"

$ws.Range("B87").Value = "x0000"
$ws.Range("C87").Value = "int sum"

$ws.Range("C88").Value = "int counter"

$ws.Range("C89").Value = "Save R1"

$ws.Range("B90").Value = "x####"
$ws.Range("C90").Value = "previous frame pointer"
$ws.Range("D90").Value = "R5"
$ws.Range("E90").Value = "current frame pointer"

$ws.Range("C91").Value = "sumOfSquares() return address"
$ws.Range("D91").Value = "R7"
$ws.Range("E91").Value = "current return address"

$ws.Range("B92").Value = "x0005"
$ws.Range("C92").Value = "arraySize parameter"

$ws.Range("B93").Value = "x5FF9"
$ws.Range("C93").Value = "array parameter"

$ws.Range("B94").Value = "x0001"
$ws.Range("C94").Value = "int total"

$ws.Range("B95").Value = "x0002"
$ws.Range("C95").Value = "array[0]"

$ws.Range("B96").Value = "x0003"
$ws.Range("C96").Value = "array[1]"

$ws.Range("B97").Value = "x0005"
$ws.Range("C97").Value = "array[2]"

$ws.Range("B98").Value = "x0000"
$ws.Range("C98").Value = "array[3]"

$ws.Range("B99").Value = "x0001"
$ws.Range("C99").Value = "array[4]"

$ws.Range("D100").Value = "R5"
$ws.Range("E100").Value = "current frame pointer"

# --- restore the active selection to C92, matching the saved view state ---
$ws.Range("C92").Select()
